# Laboratorio 7 - Doc Pruebas - Est 1
# Actualiza los datos medidos (Consumo de Datos y Tiempo de Ejecucion Real)
# de las tablas PROBING y CHAINING en la hoja "Datos Lab7", reemplazando los
# valores/formulas de muestra originales por los resultados reales de la prueba.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab7")

# --- Tabla "Carga de Catalogo PROBING (-large)" (Table1), filas 3-6 ---
$ws.Range("B3").Value = 34929.78
$ws.Range("C3").Value = 251.05

$ws.Range("B4").Value = 34929.78
$ws.Range("C4").Value = 225.16

$ws.Range("B5").Value = 34929.78
$ws.Range("C5").Value = 194.11

$ws.Range("B6").Value = 35114.089999999997
$ws.Range("C6").Value = 189.27

# --- Tabla "Carga de Catalogo CHAINING (-large)" (Table13), filas 11-14 ---
$ws.Range("B11").Value = 35114.82
$ws.Range("C11").Value = 214.91

$ws.Range("B12").Value = 35114.82
$ws.Range("C12").Value = 203.05

$ws.Range("B13").Value = 35282.559999999998
$ws.Range("C13").Value = 203.16

$ws.Range("B14").Value = 45348.3
$ws.Range("C14").Value = 206.98

# Deja el cursor donde quedo la edicion
$ws.Range("C21").Select() | Out-Null
